# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G ("K") for rows 2-16
$values = @{
    2  = 1
    3  = 2
    4  = 2
    5  = 1
    6  = 0
    7  = 0
    8  = 0
    9  = 1
    10 = 1
    11 = 3
    12 = 1
    13 = 1
    14 = 1
    15 = 2
    16 = 2
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}
